$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 9797
$ws.Range("J17").Value = 9797
$ws.Range("L17").Value = 29391
$ws.Range("N17").Value = -29727
$ws.Range("H112").Value = 1255.8823
$ws.Range("J112").Value = 1271.875
$ws.Range("L112").Value = 3815.625
$ws.Range("N112").Value = -6031.625
$ws.Range("H125").Value = 1431.2
$ws.Range("I125").Value = 1289
$ws.Range("K125").Value = 11601
$ws.Range("M125").Value = -9141
$ws.Range("H129").Value = 431479.53
$ws.Range("J129").Value = 3313.8096
$ws.Range("L129").Value = 9941.4288
$ws.Range("N129").Value = -19941.4288
$ws.Range("H137").Value = 3205.2666
$ws.Range("I137").Value = 1247.4546
$ws.Range("J137").Value = 3541.7656
$ws.Range("K137").Value = 3742.3638
$ws.Range("L137").Value = 10625.2968
$ws.Range("M137").Value = -1192.3638
$ws.Range("N137").Value = -15725.2968
$ws.Range("H138").Value = 3104.3547
$ws.Range("I138").Value = 3006.5557
$ws.Range("J138").Value = 3127.8267
$ws.Range("K138").Value = 9019.667099999999
$ws.Range("L138").Value = 9383.480100000001
$ws.Range("M138").Value = -3879.667099999999
$ws.Range("N138").Value = -19663.4801
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24422.06
$ws.Range("I32").Value = 24422.06
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 24422.06
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -24135.06
$ws.Range("N32").ClearContents()
$ws.Range("H61").Value = 1747.75
$ws.Range("I61").Value = 1515.7179
$ws.Range("K61").Value = 1515.7179
$ws.Range("M61").Value = -1303.7179
$ws.Range("H136").Value = 1747.75
$ws.Range("I136").Value = 1515.7179
$ws.Range("K136").Value = 4547.153700000001
$ws.Range("M136").Value = -1997.153700000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12521.909
$ws.Range("I31").Value = 3183.5
$ws.Range("J31").Value = 23728
$ws.Range("K31").Value = 3183.5
$ws.Range("L31").Value = 23728
$ws.Range("M31").Value = -2888.5
$ws.Range("N31").Value = -24318
$ws.Range("H34").Value = 12521.909
$ws.Range("I34").Value = 3183.5
$ws.Range("J34").Value = 23728
$ws.Range("K34").Value = 3183.5
$ws.Range("L34").Value = 23728
$ws.Range("M34").Value = -2981.5
$ws.Range("N34").Value = -24132
$ws.Range("H122").Value = 200244500
$ws.Range("J122").Value = 19914
$ws.Range("L122").Value = 59742
$ws.Range("N122").Value = -64642
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 1745.0834
$ws.Range("I8").Value = 1745.0834
$ws.Range("K8").Value = 5235.2502
$ws.Range("M8").Value = -5096.2502
$ws.Range("H131").Value = 30792.705
$ws.Range("I131").Value = 10554.9
$ws.Range("J131").Value = 34110.38
$ws.Range("K131").Value = 31664.7
$ws.Range("L131").Value = 102331.14
$ws.Range("M131").Value = -26624.7
$ws.Range("N131").Value = -112411.14
$ws.Range("H133").Value = 6265.8
$ws.Range("I133").Value = 6365
$ws.Range("J133").Value = 6199.6665
$ws.Range("K133").Value = 19095
$ws.Range("L133").Value = 18598.9995
$ws.Range("M133").Value = -14035
$ws.Range("N133").Value = -28718.9995
$ws.Range("H136").Value = 45457360
$ws.Range("I136").Value = 71430270
$ws.Range("J136").Value = 4758.25
$ws.Range("K136").Value = 214290810
$ws.Range("L136").Value = 14274.75
$ws.Range("M136").Value = -214285710
$ws.Range("N136").Value = -24474.75
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 241800.23
$ws.Range("I80").Value = 504480.5
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 504480.5
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = -503482.5
$ws.Range("N80").Value = -4996
$ws.Range("H83").Value = 241800.23
$ws.Range("I83").Value = 504480.5
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 2522402.5
$ws.Range("L83").Value = 15000
$ws.Range("M83").Value = -2517410.5
$ws.Range("N83").Value = -24984
$ws.Range("H102").Value = 1668.6471
$ws.Range("I102").Value = 1471.0834
$ws.Range("J102").Value = 2142.8
$ws.Range("K102").Value = 1471.0834
$ws.Range("L102").Value = 2142.8
$ws.Range("M102").Value = 150.9166
$ws.Range("N102").Value = -5386.8
$ws.Range("H126").Value = 13122.833
$ws.Range("I126").Value = 18434.25
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 55302.75
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = -52832.75
$ws.Range("N126").Value = -12440
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2700.8
$ws.Range("I7").Value = 2700.8
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 2700.8
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -2588.8
$ws.Range("N7").ClearContents()
$ws.Range("H40").Value = 2483.0667
$ws.Range("I40").Value = 2365.1538
$ws.Range("J40").Value = 3249.5
$ws.Range("K40").Value = 2365.1538
$ws.Range("L40").Value = 3249.5
$ws.Range("M40").Value = -2229.1538
$ws.Range("N40").Value = -3521.5
$ws.Range("H122").Value = 33129.188
$ws.Range("I122").Value = 33129.188
$ws.Range("K122").Value = 99387.56400000001
$ws.Range("M122").Value = -96937.56400000001
$ws.Range("H126").Value = 2700.8
$ws.Range("I126").Value = 2700.8
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 8102.400000000001
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -5632.400000000001
$ws.Range("N126").ClearContents()
$ws.Range("H136").Value = 2303.853
$ws.Range("I136").Value = 1762.64
$ws.Range("J136").Value = 3807.2222
$ws.Range("K136").Value = 5287.92
$ws.Range("L136").Value = 11421.6666
$ws.Range("M136").Value = -2737.92
$ws.Range("N136").Value = -16521.6666
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H124").Value = 45429
$ws.Range("J124").Value = 45429
$ws.Range("L124").Value = 45429
$ws.Range("N124").Value = -55249
$ws.Range("H126").Value = 2263022
$ws.Range("I126").Value = 2674207.8
$ws.Range("J126").Value = 1500
$ws.Range("K126").Value = 8022623.399999999
$ws.Range("L126").Value = 4500
$ws.Range("M126").Value = -8020153.399999999
$ws.Range("N126").Value = -9440
$ws.Range("H136").Value = 16097.471
$ws.Range("I136").Value = 22685.39
$ws.Range("J136").Value = 2322.7273
$ws.Range("K136").Value = 68056.17
$ws.Range("L136").Value = 6968.1819
$ws.Range("M136").Value = -65506.17
$ws.Range("N136").Value = -12068.1819
